# Update "want-to-go" counts (column F) on the "展览" (Exhibition) sheet
# and on the corresponding rows of the "全部类型" (All types) sheet, to
# match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# 展览 sheet
$sheetExhibition.Range("F3").Value  = 28
$sheetExhibition.Range("F9").Value  = 555
$sheetExhibition.Range("F13").Value = 13523
$sheetExhibition.Range("F17").Value = 5561
$sheetExhibition.Range("F18").Value = 5584
$sheetExhibition.Range("F19").Value = 58

# 全部类型 sheet (same rows of data, different row offsets)
$sheetAll.Range("F10").Value = 28
$sheetAll.Range("F31").Value = 555
$sheetAll.Range("F35").Value = 13523
$sheetAll.Range("F40").Value = 5561
$sheetAll.Range("F41").Value = 5584
$sheetAll.Range("F42").Value = 58
